$d = $word.ActiveDocument

# Helper: wdReplaceOne = 1 (replace only the first/only match), and we
# always scope the Find to the specific paragraph's Range so edits in
# one paragraph can never bleed into another paragraph that happens to
# contain matching text.

# ------------------------------------------------------------------
# Locate the stable paragraphs we need to touch, by their (unchanging)
# original text, before any of them get mutated.
# ------------------------------------------------------------------
$docentePara = $null
$resumidoPara = $null
$programaPara = $null
$avaliacaoPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -eq "3577649 - Carlos Angelo Nunes`r") {
        $docentePara = $p
    } elseif ($t -eq "1176388 - Luiz Tadeu Fernandes Eleno`r") {
        $resumidoPara = $p
    } elseif ($t -eq "1) 1a Lei da Termodinâmica2) 2a e 3a Leis da Termodinâmica3) Relações entre Propriedades Termodinâmicas4) Equilíbrio5) Equilíbrio Químico6) Soluções`r") {
        $programaPara = $p
    } elseif ($t.StartsWith("Método: ")) {
        $avaliacaoPara = $p
    }
}

# ------------------------------------------------------------------
# 2) The paragraph that used to hold "1176388 - Luiz Tadeu..." now
#    holds the short "Programa resumido" list text.
# ------------------------------------------------------------------
$resumidoPara.Range.Find.Execute("1176388 - Luiz Tadeu Fernandes Eleno", $true, $false, $false, $false, $false, $true, 1, $false, "1) 1a Lei da Termodinâmica2) 2a e 3a Leis da Termodinâmica3) Relações entre Propriedades Termodinâmicas4) Equilíbrio5) Equilíbrio Químico6) Soluções", 1) | Out-Null

# ------------------------------------------------------------------
# 3) The paragraph that used to hold the short "Programa resumido"
#    text now holds the long detailed "Programa" description (which
#    used to be the "Método" run content in "Avaliação").
# ------------------------------------------------------------------
$programaPara.Range.Find.Execute("1) 1a Lei da Termodinâmica2) 2a e 3a Leis da Termodinâmica3) Relações entre Propriedades Termodinâmicas4) Equilíbrio5) Equilíbrio Químico6) Soluções", $true, $false, $false, $false, $false, $true, 1, $false, "1- 1a Lei da Termodinâmica: sistema e vizinhança; transferência de energia; energia de um sistema; energia como uma função de estado; trabalho; sistema fechado; propriedades extensivas e intensivas; sistema aberto; entalpia; estado estacionário; capacidade térmica a volume constante; capacidade térmica a volume constante; equação de estado, gases não-ideais; expansão e compressão adiabática; entalpias de formação; variação de entalpia em reações químicas; variação de temperatura associadas à reações químicas em sistemas adiabáticos.2- 2a e 3a Leis da Termodinâmica: Entropia como função de estado; variações de entropia associadas à variações de temperatura e pressão; variações de entropia em reações químicas; terceiro princípio da termodinâmica.3- Relações entre Propriedades Termodinâmicas: As funções A e G; potencial químico; grandezas parciais molares; relações entre propriedades derivadas de U, H, A e G; gás ideal; entropia de mistura; capacidade térmica; variação de capacidade térmica; Relação P-T isoentrópica; compressão isoentrópica de sólidos.4- Equilíbrio: Condições de equilíbrio; equilíbrio de fases; variação de pressão de equilíbrio com a temperatura; equação de Clapeyron; variação da pressão de vapor de uma fase condensada com a pressão total aplicada; variação da pressão de vapor com tamanho de partícula.5- Equilíbrio Químico: atividade; equilíbrio químico; equilíbrio em fase gasosa; equilíbrio sólido-vapor; fontes de informação em DGo; diagrama de Ellingham; variação da constante de equilíbrio com a temperatura; gases dissolvidos em metais (Lei de Sievert); equilíbrio químico e temperatura adiabática de chama.6- Soluções: grandezas parciais molares relativas; entropia de mistura - solução ideal; entalpia de mistura  solução ideal; Soluções não-ideais; relação de Gibbs-Duhem; soluções regulares.", 1) | Out-Null

# ------------------------------------------------------------------
# 4) In the "Avaliação" paragraph, the content of each labeled run
#    shifts up: Método <- (old Programa long text), Critério <- (old
#    Método text), Norma de recuperação <- (old Critério text), and
#    a brand-new closing text for Norma de recuperação.
#    Do these back-to-front so a find for an "old" value never
#    accidentally matches text we *just* wrote earlier in this same
#    paragraph.
# ------------------------------------------------------------------
$avaliacaoPara.Range.Find.Execute("Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF). A nota final será calculada através da expressão:NF=(P1+2*P2)/3", $true, $false, $false, $false, $false, $true, 1, $false, "Para os alunos que obtiverem 3,0≤NF<5,0, será aplicada uma avaliação de recuperação (R), com pontuação de 0 a 10, que levará ao cálculo da média final (MF) através da seguinte expressão:MF=(NF+R)/2", 1) | Out-Null

$avaliacaoPara.Range.Find.Execute("O curso será ministrado na forma de aulas expositivas.", $true, $false, $false, $false, $false, $true, 1, $false, "Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF). A nota final será calculada através da expressão:NF=(P1+2*P2)/3", 1) | Out-Null

$avaliacaoPara.Range.Find.Execute("1- 1a Lei da Termodinâmica: sistema e vizinhança; transferência de energia; energia de um sistema; energia como uma função de estado; trabalho; sistema fechado; propriedades extensivas e intensivas; sistema aberto; entalpia; estado estacionário; capacidade térmica a volume constante; capacidade térmica a volume constante; equação de estado, gases não-ideais; expansão e compressão adiabática; entalpias de formação; variação de entalpia em reações químicas; variação de temperatura associadas à reações químicas em sistemas adiabáticos.2- 2a e 3a Leis da Termodinâmica: Entropia como função de estado; variações de entropia associadas à variações de temperatura e pressão; variações de entropia em reações químicas; terceiro princípio da termodinâmica.3- Relações entre Propriedades Termodinâmicas: As funções A e G; potencial químico; grandezas parciais molares; relações entre propriedades derivadas de U, H, A e G; gás ideal; entropia de mistura; capacidade térmica; variação de capacidade térmica; Relação P-T isoentrópica; compressão isoentrópica de sólidos.4- Equilíbrio: Condições de equilíbrio; equilíbrio de fases; variação de pressão de equilíbrio com a temperatura; equação de Clapeyron; variação da pressão de vapor de uma fase condensada com a pressão total aplicada; variação da pressão de vapor com tamanho de partícula.5- Equilíbrio Químico: atividade; equilíbrio químico; equilíbrio em fase gasosa; equilíbrio sólido-vapor; fontes de informação em DGo; diagrama de Ellingham; variação da constante de equilíbrio com a temperatura; gases dissolvidos em metais (Lei de Sievert); equilíbrio químico e temperatura adiabática de chama.6- Soluções: grandezas parciais molares relativas; entropia de mistura - solução ideal; entalpia de mistura  solução ideal; Soluções não-ideais; relação de Gibbs-Duhem; soluções regulares.", $true, $false, $false, $false, $false, $true, 1, $false, "O curso será ministrado na forma de aulas expositivas.", 1) | Out-Null

# ------------------------------------------------------------------
# 1) "Docente(s) Responsável(eis)" paragraph: append a line break and
#    a new run with the text that used to live in its own paragraph
#    ("1176388 - Luiz Tadeu Fernandes Eleno"). Do this last, since it
#    re-introduces that literal text into the document and could
#    otherwise confuse earlier (already-completed) Find/Replace steps.
# ------------------------------------------------------------------
$r = $docentePara.Range
$insertPoint = $r.End - 1
$breakRange = $d.Range($insertPoint, $insertPoint)
$breakRange.InsertBreak(6)  # wdLineBreak
$textRange = $d.Range($insertPoint + 1, $insertPoint + 1)
$textRange.InsertAfter("1176388 - Luiz Tadeu Fernandes Eleno")
